# Reports.xlsx: refresh the OFFICE_CHECK demo data on Sheet1
# (4th time using it) - replaces the existing 18-row Firm/Time/Lawyers
# Registered sample table (rows 2:19) with a freshly regenerated 50-row
# table (rows 2:51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('SFKS Law', '10s', '1'),
  @('Howse Williams', '11s', '1'),
  @('Foyen', '08s', '1'),
  @('Ritch Mueller And Nicolau', '07s', '1'),
  @('Dechert LLP', '09s', '3'),
  @('SEUM', '41s', '1'),
  @('Krogerus', '08s', '1'),
  @('Njord', '05s', '2'),
  @('AWA', '20s', '3'),
  @('Finreg 360', '06s', '1'),
  @('Taylor Wessing', '29s', '3'),
  @('VO Patents And Trademarks', '09s', '2'),
  @('MZM Legal', '23s', '1'),
  @('William Fry', '34s', '2'),
  @('Ashurst', '08s', '3'),
  @('Tavares', '05s', '1'),
  @('LEFOSSE', '06s', '1'),
  @('DaleAndLessmann', '06s', '1'),
  @('MBIP', '42s', '1'),
  @('Mijares Angoitia Cortés And Fuentes', '09s', '1'),
  @('Ropes And Gray', '14s', '3'),
  @('CFN Law', '08s', '1'),
  @('Kolster', '07s', '1'),
  @('Collas Crill', '12s', '3'),
  @('BLG Law', '12s', '1'),
  @('Gianni And Origoni', '37s', '3'),
  @('Stikeman Elliott', '01min 19s', '0'),
  @('Jones Day', '44s', '3'),
  @('BARDEHLE PAGENBERG', '06s', '2'),
  @('Viering Jentschura And Partner', '09s', '1'),
  @('Minter Ellison RuddWatts', '12s', '1'),
  @('Carey Olsen', '22s', '3'),
  @('Deacons', '13s', '1'),
  @('Panetta Consulting Group', '06s', '1'),
  @('Havel Partners', '10s', '1'),
  @('Bennett Jones', '08s', '0'),
  @('Hammarskiöld And Co', '09s', '1'),
  @('Meitar Law', '01min 09s', '1'),
  @('Consortium Legal', '12s', '3'),
  @('Kromann Reumert', '16s', '1'),
  @('Malley And Co', '03s', '1'),
  @('Werksmans', '11s', '1'),
  @('Spruson And Ferguson', '01min 13s', '3'),
  @('Drzewiecki Tomaszek', '05s', '1'),
  @('Aguayo Ecclefield And Martinez', '06s', '1'),
  @('Appleby Global', '10s', '3'),
  @('Titov', '09s', '1'),
  @('LEX Logmannsstofa', '10s', '1'),
  @('Campbells Legal', '08s', '2'),
  @('Arnold And Porter', '20s', '3')
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]

  # Column C ("Lawyers Registered") holds small digit strings ("0".."3").
  # A plain .Value assignment would be auto-coerced to a number, so the
  # value is written as a text formula and then pasted back as a value -
  # this keeps the cell a genuine text/shared-string cell (matching the
  # rest of the table) instead of a numeric one.
  $ws.Cells.Item($r, 3).Formula = '="' + $row[2] + '"'
  $ws.Cells.Item($r, 3).Copy()
  $ws.Cells.Item($r, 3).PasteSpecial(-4163)

  $r = $r + 1
}
